$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rename the "deek_linking" sheet to "deep_linking" and make it the active
#    tab (Activate() also flips workbook bookViews/activeTab + sheetView
#    tabSelected for us).
# ---------------------------------------------------------------------------
$wsDeep = $wb.Worksheets.Item("deek_linking")
$wsDeep.Name = "deep_linking"

# ---------------------------------------------------------------------------
# 2. Strip all the hyperlinks out of the deep_linking sheet and clear the
#    leftover "Hyperlink" character formatting (blue/underlined) so the
#    cells fall back to the default style.
# ---------------------------------------------------------------------------
$wsDeep.Hyperlinks.Delete() | Out-Null
$linkRange = $wsDeep.Range("A2:B19")
$linkRange.Font.Underline = $false
$linkRange.Font.ColorIndex = 0
$linkRange.ClearFormats() | Out-Null

# ---------------------------------------------------------------------------
# 3. Update headers + "record type" column text for the new table layout.
# ---------------------------------------------------------------------------
$wsDeep.Cells.Item(1, 2).Value = "record type"
$wsDeep.Cells.Item(1, 3).Value = "Run Mode"

$recordTypes = @{
    2  = "For All search page "
    3  = "For Article search page "
    4  = "For Patent search page "
    5  = "For People search page "
    6  = "For Post search page "
    7  = "For All search page "
    8  = "For All search page "
    9  = "For All search page "
    10 = "For Article  search page "
    11 = "For Article search page "
    12 = "For Article search page "
    13 = "For Patents search page "
    14 = "For patents search page "
    15 = "For patents search page "
    16 = "For People search page "
    17 = "For Posts search page "
    18 = "For Posts search page "
    19 = "For Posts search page "
}
foreach ($row in $recordTypes.Keys) {
    $wsDeep.Cells.Item($row, 2).Value = $recordTypes[$row]
}

# ---------------------------------------------------------------------------
# 4. Widen/adjust the deep_linking columns and add the new "Run Mode" column
#    width, matching the refreshed layout.
# ---------------------------------------------------------------------------
$wsDeep.Columns.Item(1).ColumnWidth = 89.28515625
$wsDeep.Columns.Item(2).ColumnWidth = 24
$wsDeep.Columns.Item(3).ColumnWidth = 10

# ---------------------------------------------------------------------------
# 5. Selection / scroll position on the deep_linking sheet.
# ---------------------------------------------------------------------------
$wsDeep.Activate()
$wsDeep.Range("A19").Select() | Out-Null

# ---------------------------------------------------------------------------
# 6. Fix up row 107 on the "Test Cases" sheet (new test case row content).
# ---------------------------------------------------------------------------
$wsCases = $wb.Worksheets.Item("Test Cases")
$wsCases.Cells.Item(107, 1).Value = "Search127"
$wsCases.Cells.Item(107, 2).Value = "OPQA-2801|OPQA-2802|OPQA-2803|OPQA-2804|OPQA-2805|OPQA-2808|OPQA-2809|OPQA-2810|OPQA-2811|OPQA-2812|OPQA-2813|OPQA-2814|OPQA-2815|OPQA-2816|OPQA-2817|OPQA-2818|OPQA-2819|OPQA-2820"
$wsCases.Cells.Item(107, 3).Value = "Verify that Deeplinking is working for Search result page using steam account"

$wsCases.Activate()
$excel.ActiveWindow.ScrollRow = 103
$excel.ActiveWindow.ScrollColumn = 1
$wsCases.Range("C112").Select() | Out-Null

# ---------------------------------------------------------------------------
# 7. Selection change on the "Test Case Steps" sheet.
# ---------------------------------------------------------------------------
$wsSteps = $wb.Worksheets.Item("Test Case Steps")
$wsSteps.Activate()
$wsSteps.Range("C14").Select() | Out-Null

# ---------------------------------------------------------------------------
# 8. Re-activate deep_linking last so it ends up the active tab (matches
#    workbookView activeTab="1" / sheetView tabSelected on deep_linking).
# ---------------------------------------------------------------------------
$wsDeep.Activate()
$wsDeep.Range("A19").Select() | Out-Null

Write-Output "edit complete"
